# Transition rule 5 and 10 mile radius updates
$wb = $excel.ActiveWorkbook

$wsMeans = $wb.Worksheets.Item("Means")
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# --- Sheet "Means" ---
# Header row
$wsMeans.Range("F1").Value = "Within 5 miles of HFC production facility"
$wsMeans.Range("G1").Value = "Within 10 miles of HFC production facility"

# Data rows: Variable, National Avg, State Avg, 1mi, 3mi, 5mi(new), 10mi(new)
$meansData = @(
    @("% White", 72, 83, 55, 40, 49, 26),
    @("% Black or African American ", 13, 9.4, 4.5, 23, 28, 64),
    @("% Other", 15, 7.3, 41, 37, 23, 9.2),
    @("% Hispanic", 18, 6.9, 55, 53, 52, 17),
    @("Median Income [1,000 2019$]", 71, 62, 39, 40, 43, 45),
    @("% Below Poverty Line", 7.3, 7, 11, 14, 13, 11),
    @("% Below Half the Poverty Line", 5.8, 6, 5.5, 12, 11, 11),
    @("Total Cancer Risk (per million)", 29, 23, 30, 30, 30, 30),
    @("Total Respiratory (hazard quotient)", 0.37, 0.3, 0.4, 0.4, 0.38, 0.38)
)

$rowIdx = 2
foreach ($row in $meansData) {
    $wsMeans.Cells.Item($rowIdx, 1).Value = $row[0]
    $wsMeans.Cells.Item($rowIdx, 2).Value = $row[1]
    $wsMeans.Cells.Item($rowIdx, 3).Value = $row[2]
    $wsMeans.Cells.Item($rowIdx, 4).Value = $row[3]
    $wsMeans.Cells.Item($rowIdx, 5).Value = $row[4]
    $wsMeans.Cells.Item($rowIdx, 6).Value = $row[5]
    $wsMeans.Cells.Item($rowIdx, 7).Value = $row[6]
    $rowIdx++
}

# --- Sheet "Standard Deviations" ---
# Header row
$wsSD.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$wsSD.Range("G1").Value = "Within 10 mile of HFC production facility SD"

$sdData = @(
    @("% White", 27, 23, 24, 26, 28, 30),
    @("% Black or African American ", 23, 19, 7.2, 26, 29, 37),
    @("% Other", 16, 8.9, 21, 19, 20, 13),
    @("% Hispanic", 22, 11, 25, 22, 26, 23),
    @("Median Income [1,000 2019$]", 37, 25, 5.3, 15, 15, 19),
    @("% Below Poverty Line", 8.7, 8.4, 6, 9.5, 9.6, 10),
    @("% Below Half the Poverty Line", 7.8, 8.4, 4.8, 11, 10, 11),
    @("Total Cancer Risk (per million)", 10, 4.7, 0, 0, 0, 1.7),
    @("Total Respiratory (hazard quotient)", 0.14, 0.058, 0.000000000000000032, 0.025, 0.042, 0.04)
)

$rowIdx = 2
foreach ($row in $sdData) {
    $wsSD.Cells.Item($rowIdx, 1).Value = $row[0]
    $wsSD.Cells.Item($rowIdx, 2).Value = $row[1]
    $wsSD.Cells.Item($rowIdx, 3).Value = $row[2]
    $wsSD.Cells.Item($rowIdx, 4).Value = $row[3]
    $wsSD.Cells.Item($rowIdx, 5).Value = $row[4]
    $wsSD.Cells.Item($rowIdx, 6).Value = $row[5]
    $wsSD.Cells.Item($rowIdx, 7).Value = $row[6]
    $rowIdx++
}
